$d = $word.ActiveDocument

# The document contains duplicate "featured image" pictures: one legitimate
# inline picture right after the title (kept), and two duplicate copies of
# the very same picture that were accidentally inserted right after the
# "Introducción" heading and right after the "Presentación del proyecto"
# heading. Those two duplicate picture paragraphs must be removed entirely
# (the paragraph, its formatting and the drawing run inside it).

# Walk the paragraphs from the end to the start so that deleting a paragraph
# doesn't shift the indices of paragraphs we still need to inspect.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)

    if ($p.Range.InlineShapes.Count -gt 0) {
        # Look at the paragraph immediately preceding this picture paragraph;
        # duplicated captures always sit directly under one of these two
        # headings, while the legitimate featured image sits right under the
        # title (Heading1), which we must keep untouched.
        $prevIndex = $i - 1
        if ($prevIndex -ge 1) {
            $prevText = $d.Paragraphs.Item($prevIndex).Range.Text.Trim()
            if ($prevText -eq "Introducción" -or $prevText -eq "Presentación del proyecto") {
                $p.Range.Delete()
            }
        }
    }
}
